$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B16: new test scenario text (wraps onto taller row 17)
$ws.Range("B16").Value = "platform compatibility testing, graph flexibility"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("B16").VerticalAlignment = -4108

# Update E16: new test steps text
$ws.Range("E16").Value = "1. install app to different version of Iphone 2.check UI and function"
$ws.Range("E16").WrapText = $true
$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E16").VerticalAlignment = -4108

# Row 17 (merged continuation row) gets a taller custom height to fit the wrapped text
$ws.Rows("17").RowHeight = 27

# Match B17/E17 (merged continuation cells) to the same wrap style
$ws.Range("B17").WrapText = $true
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("B17").VerticalAlignment = -4108
$ws.Range("E17").WrapText = $true
$ws.Range("E17").HorizontalAlignment = -4108
$ws.Range("E17").VerticalAlignment = -4108

# Sheet view changes: zoom + selection
$excel.ActiveWindow.Zoom = 125
$ws.Range("I13").Select()

# Page setup
$ws.PageSetup.Orientation = 1
